$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for the 4 new rows (25-28): Chest Reopening Concept 3 vertices, printed on all printers
$dateStr = "27-06-2018"

$rows = @(
    @{ Row = 25; Part = "Chest Reopening Concept 3 Top Front Verticies"; Material = "PLA" },
    @{ Row = 26; Part = "Chest Reopening Concept 3 Top Back Verticies";  Material = "Polylite" },
    @{ Row = 27; Part = "Chest Reopening Concept 3 Bot Front Verticies"; Material = "Polylite" },
    @{ Row = 28; Part = "Chest Reopening Concept 3 Bot Back Verticies";  Material = "Polylite" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $dateStr
    $ws.Cells.Item($row, 2).Value = $dateStr
    $ws.Cells.Item($row, 3).Value = $r.Part
    $ws.Cells.Item($row, 4).Value = 2
    $ws.Cells.Item($row, 5).Value = $r.Material
    $ws.Cells.Item($row, 6).Value = 2
    $ws.Cells.Item($row, 7).Value = 20
    $ws.Cells.Item($row, 8).Value = 0.2
    $ws.Cells.Item($row, 9).Value = "NA"
}

$ws.Columns("C:C").ColumnWidth = 37.17

$ws.Range("J28:K28").Select() | Out-Null
